$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.168.46"
$ws.Range("E2").Value = "  -6.20%  "
$ws.Range("D3").Value = "3.283.08"
$ws.Range("E3").Value = "  -5.49%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'557.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").Value = "'127.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.275.56"
$ws.Range("E8").Value = "  -5.68%  "
$ws.Range("D9").Value = "'0.472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "'7.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.19%  "
$ws.Range("E11").Value = "  -5.08%  "
$ws.Range("E12").Value = "  -4.59%  "
$ws.Range("D13").Value = "3.832.17"
$ws.Range("E13").Value = "  -5.70%  "
$ws.Range("D14").Value = "'0.120"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "3.278.16"
$ws.Range("E15").Value = "  -5.57%  "
$ws.Range("E16").Value = "  -5.86%  "
$ws.Range("D17").Value = "'24.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "60.088.54"
$ws.Range("E18").Value = "  -6.19%  "
$ws.Range("D19").Value = "'5.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "'13.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "'9.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.52%  "
$ws.Range("D22").Value = "'351.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.45%  "
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "3.400.47"
$ws.Range("E25").Value = "  -5.85%  "
$ws.Range("E26").Value = "  -7.64%  "
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "'0.984"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").Value = "'7.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "'7.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("E32").Value = "  -6.26%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("D35").Value = "3.300.26"
$ws.Range("E35").Value = "  -5.70%  "
$ws.Range("D36").Value = "'22.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "'5.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").Value = "'6.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "'157.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("D41").Value = "'0.0745"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("D42").Value = "'0.996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "'40.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  -7.64%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  -4.94%  "
$ws.Range("D48").Value = "'22.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "'0.855"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.23%  "
$ws.Range("D51").Value = "'21.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.29%  "
